$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 21.1
$ws.Range("I11").Value = 21.1
$ws.Range("K11").Value = 21.1
$ws.Range("M11").Value = 118.9
$ws.Range("H33").Value = 392.65714
$ws.Range("I33").Value = 204.54839
$ws.Range("K33").Value = 204.54839
$ws.Range("M33").Value = 24.45160999999999
$ws.Range("H40").Value = 2741.6667
$ws.Range("I40").Value = 0
$ws.Range("J40").Value = 2741.6667
$ws.Range("K40").Value = 0
$ws.Range("L40").ClearContents()
$ws.Range("M40").Value = 2741.6667
$ws.Range("N40").Value = -3091.6667
$ws.Range("H46").Value = 27825
$ws.Range("I46").Value = 650
$ws.Range("J46").Value = 55000
$ws.Range("K46").Value = 1950
$ws.Range("L46").Value = 165000
$ws.Range("M46").Value = -1831
$ws.Range("N46").Value = -165238
$ws.Range("H60").Value = 27825
$ws.Range("I60").Value = 650
$ws.Range("J60").Value = 55000
$ws.Range("K60").Value = 1950
$ws.Range("L60").Value = 165000
$ws.Range("M60").Value = -1466
$ws.Range("N60").Value = -165968
$ws.Range("H87").Value = 29500
$ws.Range("J87").Value = 29500
$ws.Range("L87").Value = 29500
$ws.Range("N87").Value = -31996
$ws.Range("H90").Value = 29500
$ws.Range("J90").Value = 29500
$ws.Range("L90").Value = 88500
$ws.Range("N90").Value = -100980
$ws.Range("H135").Value = 333334700
$ws.Range("I135").Value = 2032
$ws.Range("J135").Value = 500001000
$ws.Range("K135").Value = 18288
$ws.Range("L135").Value = 4500009000
$ws.Range("M135").Value = -15753
$ws.Range("N135").Value = -4500014070

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 7710.25
$ws.Range("I32").Value = 3815.3767
$ws.Range("J32").Value = 25626.666
$ws.Range("K32").Value = 3815.3767
$ws.Range("L32").Value = 25626.666
$ws.Range("M32").Value = -3528.3767
$ws.Range("N32").Value = -26200.666
$ws.Range("H45").Value = 10103357
$ws.Range("I45").Value = 15153518
$ws.Range("J45").Value = 3033.3333
$ws.Range("K45").Value = 15153518
$ws.Range("L45").Value = 3033.3333
$ws.Range("M45").Value = -15153141
$ws.Range("N45").Value = -3787.3333
$ws.Range("H88").Value = 1920
$ws.Range("I88").Value = 2100
$ws.Range("J88").Value = 1800
$ws.Range("K88").Value = 2100
$ws.Range("L88").Value = 1800
$ws.Range("M88").Value = -1694
$ws.Range("N88").Value = -2612
$ws.Range("H91").Value = 1920
$ws.Range("I91").Value = 2100
$ws.Range("J91").Value = 1800
$ws.Range("K91").Value = 2100
$ws.Range("L91").Value = 1800
$ws.Range("M91").Value = -696
$ws.Range("N91").Value = -4608
$ws.Range("H98").Value = 25750
$ws.Range("J98").Value = 25750
$ws.Range("L98").Value = 25750
$ws.Range("N98").Value = -31740
$ws.Range("H122").Value = 1685.4
$ws.Range("I122").Value = 1490.579
$ws.Range("J122").Value = 2302.3333
$ws.Range("K122").Value = 4471.737
$ws.Range("L122").Value = 6906.999899999999
$ws.Range("M122").Value = -2021.737
$ws.Range("N122").Value = -11806.9999
$ws.Range("H133").Value = 37000
$ws.Range("J133").Value = 37000
$ws.Range("L133").Value = 37000
$ws.Range("N133").Value = -42060

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H75").Value = 14474.75
$ws.Range("I75").Value = 9966.333000000001
$ws.Range("J75").Value = 28000
$ws.Range("K75").Value = 9966.333000000001
$ws.Range("L75").Value = 28000
$ws.Range("M75").Value = -9030.333000000001
$ws.Range("N75").Value = -29872
$ws.Range("H78").Value = 14474.75
$ws.Range("I78").Value = 9966.333000000001
$ws.Range("J78").Value = 28000
$ws.Range("K78").Value = 29898.999
$ws.Range("L78").Value = 84000
$ws.Range("M78").Value = -25218.999
$ws.Range("N78").Value = -93360
$ws.Range("H86").Value = 1766.5
$ws.Range("I86").Value = 1561.3334
$ws.Range("J86").Value = 2053.7334
$ws.Range("K86").Value = 1561.3334
$ws.Range("L86").Value = 2053.7334
$ws.Range("M86").Value = -438.3334
$ws.Range("N86").Value = -4299.7334
$ws.Range("H89").Value = 1766.5
$ws.Range("I89").Value = 1561.3334
$ws.Range("J89").Value = 2053.7334
$ws.Range("K89").Value = 7806.666999999999
$ws.Range("L89").Value = 10268.667
$ws.Range("M89").Value = -2190.666999999999
$ws.Range("N89").Value = -21500.667
$ws.Range("H135").Value = 36666.668
$ws.Range("J135").Value = 36666.668
$ws.Range("L135").Value = 36666.668
$ws.Range("N135").Value = -46806.668

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H28").Value = 29757.666
$ws.Range("J28").Value = 29757.666
$ws.Range("L28").Value = 29757.666
$ws.Range("N28").Value = -30247.666
$ws.Range("H86").Value = 12641.25
$ws.Range("I86").Value = 4813.5713
$ws.Range("K86").Value = 4813.5713
$ws.Range("M86").Value = -3690.5713
$ws.Range("H89").Value = 12641.25
$ws.Range("I89").Value = 4813.5713
$ws.Range("K89").Value = 24067.8565
$ws.Range("M89").Value = -18451.8565
$ws.Range("H141").Value = 88909.7
$ws.Range("J141").Value = 88909.7
$ws.Range("L141").Value = 88909.7
$ws.Range("N141").Value = -99269.7

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 4682.9375
$ws.Range("I3").Value = 2503.375
$ws.Range("J3").Value = 6862.5
$ws.Range("K3").Value = 7510.125
$ws.Range("L3").Value = 20587.5
$ws.Range("M3").Value = -7398.125
$ws.Range("N3").Value = -20811.5
$ws.Range("H23").Value = 166.56522
$ws.Range("I23").Value = 65.818184
$ws.Range("J23").Value = 258.91666
$ws.Range("K23").Value = 197.454552
$ws.Range("L23").Value = 776.7499799999999
$ws.Range("M23").Value = 37.54544799999999
$ws.Range("N23").Value = -1246.74998
$ws.Range("H62").Value = 4033.923
$ws.Range("I62").Value = 2000
$ws.Range("J62").Value = 4203.4165
$ws.Range("K62").Value = 6000
$ws.Range("L62").Value = 12610.2495
$ws.Range("M62").Value = -5314
$ws.Range("N62").Value = -13982.2495
$ws.Range("H65").Value = 4033.923
$ws.Range("I65").Value = 2000
$ws.Range("J65").Value = 4203.4165
$ws.Range("K65").Value = 18000
$ws.Range("L65").Value = 37830.7485
$ws.Range("M65").Value = -14568
$ws.Range("N65").Value = -44694.7485
$ws.Range("H98").Value = 3480
$ws.Range("I98").Value = 0
$ws.Range("J98").Value = 3480
$ws.Range("K98").Value = 0
$ws.Range("L98").ClearContents()
$ws.Range("M98").Value = 10440
$ws.Range("N98").Value = -13436
$ws.Range("H113").Value = 1010647.5
$ws.Range("I113").Value = 2755415.8
$ws.Range("J113").Value = 518.5789
$ws.Range("K113").Value = 8266247.399999999
$ws.Range("L113").Value = 1555.7367
$ws.Range("M113").Value = -8264077.399999999
$ws.Range("N113").Value = -5895.736699999999
$ws.Range("H114").Value = 19091438
$ws.Range("I114").Value = 16667014
$ws.Range("J114").Value = 22728072
$ws.Range("K114").Value = 50001042
$ws.Range("L114").Value = 68184216
$ws.Range("M114").Value = -49997788
$ws.Range("N114").Value = -68190724
$ws.Range("H131").Value = 940
$ws.Range("J131").Value = 946.46466
$ws.Range("L131").Value = 2839.39398
$ws.Range("N131").Value = -12919.39398
$ws.Range("H136").Value = 3870
$ws.Range("I136").Value = 1242.8572
$ws.Range("J136").Value = 10000
$ws.Range("K136").Value = 3728.5716
$ws.Range("L136").Value = 30000
$ws.Range("M136").Value = 1371.4284
$ws.Range("N136").Value = -40200

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H4").Value = 3698.8
$ws.Range("I4").Value = 3250
$ws.Range("K4").Value = 3250
$ws.Range("M4").Value = -3138
$ws.Range("H100").Value = 29788.75
$ws.Range("J100").Value = 29788.75
$ws.Range("L100").Value = 29788.75
$ws.Range("N100").Value = -31952.75
$ws.Range("H111").Value = 22000
$ws.Range("J111").Value = 22000
$ws.Range("L111").Value = 22000
$ws.Range("N111").Value = -28134
